$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Preserve the existing text formatting for the Price/Volume columns so that
# numeric-looking strings (e.g. "27.079.25", "1.001", "0.000008510") are not
# auto-converted to numbers by Excel when assigned via .Value.
$ws.Range("D2:E51").NumberFormat = "@"

$ws.Range("D2").Value = '27.079.25'
$ws.Range("E2").Value = '  +0.17%  '
$ws.Range("D3").Value = '1.890.45'
$ws.Range("E3").Value = '  +0.95%  '
$ws.Range("D4").Value = '1.001'
$ws.Range("E4").Value = '  +0.37%  '
$ws.Range("D5").Value = '306.73'
$ws.Range("E5").Value = '  +0.37%  '
$ws.Range("D6").Value = '1.001'
$ws.Range("E6").Value = '  +0.26%  '
$ws.Range("D7").Value = '0.5178'
$ws.Range("E7").Value = '  +2.41%  '
$ws.Range("D8").Value = '0.3753'
$ws.Range("E8").Value = '  +2.56%  '
$ws.Range("D9").Value = '0.07215'
$ws.Range("E9").Value = '  +0.43%  '
$ws.Range("D10").Value = '21.10'
$ws.Range("E10").Value = '  +1.57%  '
$ws.Range("D11").Value = '0.9005'
$ws.Range("E11").Value = '  +0.58%  '
$ws.Range("D12").Value = '1.928.46'
$ws.Range("E12").Value = '  +2.95%  '
$ws.Range("D13").Value = '0.07653'
$ws.Range("E13").Value = '  +1.91%  '
$ws.Range("D14").Value = '94.19'
$ws.Range("E14").Value = '  -1.22%  '
$ws.Range("D15").Value = '5.232'
$ws.Range("E15").Value = '  -0.23%  '
$ws.Range("D16").Value = '1.002'
$ws.Range("E16").Value = '  +0.38%  '
$ws.Range("D17").Value = '0.000008510'
$ws.Range("E17").Value = '  -0.37%  '
$ws.Range("D18").Value = '14.42'
$ws.Range("E18").Value = '  +1.23%  '
$ws.Range("D19").Value = '1.000'
$ws.Range("E19").Value = '  -0.10%  '
$ws.Range("D20").Value = '27.144.24'
$ws.Range("E20").Value = '  +0.39%  '
$ws.Range("D21").Value = '5.059'
$ws.Range("E21").Value = '  +0.64%  '
$ws.Range("D22").Value = '2.146.82'
$ws.Range("E22").Value = '  +2.40%  '
$ws.Range("D23").Value = '10.59'
$ws.Range("E23").Value = '  +1.74%  '
$ws.Range("D24").Value = '6.381'
$ws.Range("E24").Value = '  -0.76%  '
$ws.Range("D25").Value = '2.294'
$ws.Range("E25").Value = '  +9.97%  '
$ws.Range("D26").Value = '145.50'
$ws.Range("E26").Value = '  -1.83%  '
$ws.Range("B27").Value = 'Toncoin'
$ws.Range("C27").Value = 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
$ws.Range("D27").Value = '1.734'
$ws.Range("E27").Value = '  -2.37%  '
$ws.Range("B28").Value = 'EthereumClassic'
$ws.Range("C28").Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$ws.Range("D28").Value = '18.06'
$ws.Range("E28").Value = '  +0.80%  '
$ws.Range("D29").Value = '114.31'
$ws.Range("E29").Value = '  +0.81%  '
$ws.Range("D30").Value = '4.987'
$ws.Range("E30").Value = '  +6.38%  '
$ws.Range("D31").Value = '4.783'
$ws.Range("E31").Value = '  +1.76%  '
$ws.Range("D32").Value = '0.09200'
$ws.Range("E32").Value = '  +0.21%  '
$ws.Range("D33").Value = '0.05051'
$ws.Range("E33").Value = '  -1.79%  '
$ws.Range("D34").Value = '1.238'
$ws.Range("E34").Value = '  +7.01%  '
$ws.Range("D35").Value = '0.7755'
$ws.Range("E35").Value = '  +2.66%  '
$ws.Range("D36").Value = '2.971'
$ws.Range("E36").Value = '  +0.21%  '
$ws.Range("D37").Value = '3.283'
$ws.Range("E37").Value = '  +2.47%  '
$ws.Range("D38").Value = '2.594'
$ws.Range("E38").Value = '  +0.16%  '
$ws.Range("D39").Value = '0.5610'
$ws.Range("E39").Value = '  -0.21%  '
$ws.Range("D40").Value = '0.01988'
$ws.Range("E40").Value = '  -0.73%  '
$ws.Range("E41").Value = '  +0.17%  '
$ws.Range("D42").Value = '9.052'
$ws.Range("E42").Value = '  +5.35%  '
$ws.Range("D43").Value = '119.59'
$ws.Range("E43").Value = '  +2.86%  '
$ws.Range("D44").Value = '6.625'
$ws.Range("E44").Value = '  +0.53%  '
$ws.Range("D45").Value = '0.1509'
$ws.Range("E45").Value = '  +2.32%  '
$ws.Range("D46").Value = '0.4834'
$ws.Range("E46").Value = '  +2.26%  '
$ws.Range("B47").Value = 'EnergySwap'
$ws.Range("C47").Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range("D47").Value = '10.20'
$ws.Range("E47").Value = '  +1.40%  '
$ws.Range("B48").Value = 'PaxDollar'
$ws.Range("C48").Value = 'https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp'
$ws.Range("D48").Value = '1.0000'
$ws.Range("E48").Value = '  +0.14%  '
$ws.Range("D49").Value = '1.594'
$ws.Range("E49").Value = '  +1.80%  '
$ws.Range("D50").Value = '37.52'
$ws.Range("E50").Value = '  +1.75%  '
$ws.Range("D51").Value = '63.93'
$ws.Range("E51").Value = '  +1.03%  '
